$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for numeric-looking price values so they are not
# coerced into floating point numbers (source data is plain text).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '43.090.37'
$ws.Range("E2").Value = '  +0.16%  '
$ws.Range("D3").Value = '2.308.97'
$ws.Range("E3").Value = '  +0.23%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '300.49'
$ws.Range("E5").Value = '  -0.45%  '
$ws.Range("D6").Value = '98.07'
$ws.Range("E6").Value = '  -0.58%  '
$ws.Range("D7").Value = '0.514'
$ws.Range("E7").Value = '  -2.19%  '
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("D9").Value = '0.508'
$ws.Range("D10").Value = '36.05'
$ws.Range("E10").Value = '  +0.97%  '
$ws.Range("D11").Value = '0.0791'
$ws.Range("E11").Value = '  +0.14%  '
$ws.Range("D12").Value = '18.23'
$ws.Range("E12").Value = '  +1.70%  '
$ws.Range("D13").Value = '0.119'
$ws.Range("E13").Value = '  +1.97%  '
$ws.Range("D14").Value = '6.82'
$ws.Range("D15").Value = '2.665.33'
$ws.Range("E15").Value = '  +0.11%  '
$ws.Range("D16").Value = '2.298.92'
$ws.Range("E16").Value = '  +1.40%  '
$ws.Range("D17").Value = '0.782'
$ws.Range("E17").Value = '  -0.89%  '
$ws.Range("D18").Value = '43.004.50'
$ws.Range("E18").Value = '  +0.23%  '
$ws.Range("D19").Value = '12.76'
$ws.Range("E19").Value = '  -4.73%  '
$ws.Range("D20").Value = '0.0₃0906'
$ws.Range("E20").Value = '  -0.43%  '
$ws.Range("D21").Value = '6.06'
$ws.Range("E21").Value = '  -1.93%  '
$ws.Range("E22").Value = '  -0.28%  '
$ws.Range("D23").Value = '241.05'
$ws.Range("E23").Value = '  +0.62%  '
$ws.Range("D24").Value = '2.15'
$ws.Range("E24").Value = '  -0.51%  '
$ws.Range("E25").Value = '  +0.10%  '
$ws.Range("D26").Value = '2.44'
$ws.Range("E26").Value = '  +0.07%  '
$ws.Range("E27").Value = '  +0.09%  '
$ws.Range("D28").Value = '25.61'
$ws.Range("E28").Value = '  +3.35%  '
$ws.Range("D29").Value = '165.88'
$ws.Range("E29").Value = '  -1.33%  '
$ws.Range("E30").Value = '  -0.53%  '
$ws.Range("D31").Value = '9.09'
$ws.Range("E31").Value = '  -0.53%  '
$ws.Range("D32").Value = '33.27'
$ws.Range("E32").Value = '  -0.24%  '
$ws.Range("D33").Value = '4.98'
$ws.Range("E33").Value = '  +2.94%  '
$ws.Range("D34").Value = '0.999'
$ws.Range("E34").Value = '  +0.05%  '
$ws.Range("E35").Value = '  -3.06%  '
$ws.Range("D36").Value = '17.10'
$ws.Range("E36").Value = '  -6.21%  '
$ws.Range("E37").Value = '  -1.07%  '
$ws.Range("D38").Value = '0.0688'
$ws.Range("E38").Value = '  -0.35%  '
$ws.Range("E39").Value = '  -0.57%  '
$ws.Range("E40").Value = '  -1.38%  '
$ws.Range("E41").Value = '  -0.31%  '
$ws.Range("E42").Value = '  -1.66%  '
$ws.Range("D43").Value = '2.010.36'
$ws.Range("E43").Value = '  +0.61%  '
$ws.Range("D44").Value = '0.0282'
$ws.Range("E44").Value = '  -2.00%  '
$ws.Range("D45").Value = '2.16'
$ws.Range("E45").Value = '  +1.77%  '
$ws.Range("E46").Value = '  +0.67%  '
$ws.Range("D47").Value = '17.55'
$ws.Range("E47").Value = '  +0.44%  '
$ws.Range("D48").Value = '2.81'
$ws.Range("E48").Value = '  -0.66%  '
$ws.Range("D49").Value = '53.93'
$ws.Range("E49").Value = '  -1.41%  '
$ws.Range("D50").Value = '2.531.35'
$ws.Range("E50").Value = '  +0.06%  '
$ws.Range("D51").Value = '72.57'
$ws.Range("E51").Value = '  -1.09%  '
